$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-30 Friday", 2) | Out-Null
$d.Content.Find.Execute("221÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "735÷4=", 2) | Out-Null
$d.Content.Find.Execute("176÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "675÷4=", 2) | Out-Null
$d.Content.Find.Execute("490÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "955÷5=", 2) | Out-Null
$d.Content.Find.Execute("665÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "963÷5=", 2) | Out-Null
$d.Content.Find.Execute("768÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "354÷9=", 2) | Out-Null
$d.Content.Find.Execute("728÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "722÷7=", 2) | Out-Null
$d.Content.Find.Execute("547÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "615÷9=", 2) | Out-Null
$d.Content.Find.Execute("608÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "352÷7=", 2) | Out-Null
$d.Content.Find.Execute("937÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "643÷7=", 2) | Out-Null
$d.Content.Find.Execute("692÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "855÷6=", 2) | Out-Null
$d.Content.Find.Execute("783÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "256÷4=", 2) | Out-Null
$d.Content.Find.Execute("754÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "914÷4=", 2) | Out-Null
$d.Content.Find.Execute("230÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "341÷9=", 2) | Out-Null
$d.Content.Find.Execute("407÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "504÷4=", 2) | Out-Null
$d.Content.Find.Execute("186÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "133÷3=", 2) | Out-Null
$d.Content.Find.Execute("863÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "639÷5=", 2) | Out-Null
$d.Content.Find.Execute("939÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "523÷8=", 2) | Out-Null
$d.Content.Find.Execute("892÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "731÷4=", 2) | Out-Null
$d.Content.Find.Execute("982÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷8=", 2) | Out-Null
$d.Content.Find.Execute("228÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "211÷6=", 2) | Out-Null
$d.Content.Find.Execute("209÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "323÷3=", 2) | Out-Null
$d.Content.Find.Execute("259÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "329÷2=", 2) | Out-Null
$d.Content.Find.Execute("456÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "435÷9=", 2) | Out-Null
$d.Content.Find.Execute("532÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "940÷7=", 2) | Out-Null
$d.Content.Find.Execute("256÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "284÷3=", 2) | Out-Null
